# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 / 演出 / 全部类型 sheets, per the refreshed crawl output.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        3  = 3133
        7  = 2623
        10 = 27
        14 = 9867
        17 = 24
        18 = 7798
        19 = 12383
        28 = 2824
        29 = 1334
        30 = 214
        33 = 4597
        34 = 1279
    }
    "演出" = @{
        6 = 1194
    }
    "全部类型" = @{
        4  = 3133
        9  = 2623
        13 = 27
        17 = 9867
        20 = 24
        21 = 7798
        22 = 12383
        34 = 214
        37 = 4597
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
